$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.09798757153841207
$ws.Range("C2").Value = 0.7301515166201069
$ws.Range("B3").Value = 0.1047583938750073
$ws.Range("C3").Value = 1.102359249518327
